# Regenerate the "K" column (column G) values for save_data sheet.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# Only column G ("K") values change; everything else stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 3
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    19 = 3
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 2
    38 = 2
    39 = 4
    40 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 1
    48 = 1
    49 = 0
    50 = 0
    51 = 1
    52 = 0
    53 = 0
    54 = 2
    55 = 2
    56 = 1
    57 = 2
    58 = 3
    59 = 1
    60 = 1
    62 = 1
    63 = 1
    64 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
